$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 13-17 (data now stops at row 12)
$ws.Range("A13:H17").Delete(-4162)

# Add header values for new columns I and J, copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8

# Update B2:H12 values and add I2:J12 values
$ws.Range("B2").Value = 0.9419999999999999
$ws.Range("C2").Value = 0.873
$ws.Range("D2").Value = -2.993
$ws.Range("E2").Value = 1.018
$ws.Range("F2").Value = 7.254
$ws.Range("G2").Value = 0.141
$ws.Range("H2").Value = 0.127
$ws.Range("I2").Value = 412
$ws.Range("J2").Value = 1919

$ws.Range("B3").Value = 0.97
$ws.Range("C3").Value = 0.921
$ws.Range("D3").Value = -3.003
$ws.Range("E3").Value = 1.005
$ws.Range("F3").Value = 5.634
$ws.Range("G3").Value = 0.08500000000000001
$ws.Range("H3").Value = 0.079
$ws.Range("I3").Value = 378
$ws.Range("J3").Value = 1759

$ws.Range("B4").Value = 0.971
$ws.Range("C4").Value = 0.895
$ws.Range("D4").Value = -3.015
$ws.Range("E4").Value = 1.034
$ws.Range("F4").Value = 4.181
$ws.Range("G4").Value = 0.115
$ws.Range("H4").Value = 0.105
$ws.Range("I4").Value = 152
$ws.Range("J4").Value = 695

$ws.Range("B5").Value = 1.006
$ws.Range("C5").Value = 0.848
$ws.Range("D5").Value = -3.023
$ws.Range("E5").Value = 1.047
$ws.Range("F5").Value = 3.145
$ws.Range("G5").Value = 0.161
$ws.Range("H5").Value = 0.152
$ws.Range("I5").Value = 178
$ws.Range("J5").Value = 815

$ws.Range("B6").Value = 1.043
$ws.Range("C6").Value = 0.873
$ws.Range("D6").Value = -3.06
$ws.Range("E6").Value = 1.053
$ws.Range("F6").Value = 2.218
$ws.Range("G6").Value = 0.156
$ws.Range("H6").Value = 0.127
$ws.Range("I6").Value = 206
$ws.Range("J6").Value = 942

$ws.Range("B7").Value = 1.038
$ws.Range("C7").Value = 0.918
$ws.Range("D7").Value = -3.019
$ws.Range("E7").Value = 1.022
$ws.Range("F7").Value = 1.335
$ws.Range("G7").Value = 0.095
$ws.Range("H7").Value = 0.082
$ws.Range("I7").Value = 209
$ws.Range("J7").Value = 957

$ws.Range("B8").Value = 1.026
$ws.Range("C8").Value = 0.946
$ws.Range("D8").Value = -2.976
$ws.Range("E8").Value = 0.999
$ws.Range("F8").Value = 0.576
$ws.Range("G8").Value = 0.065
$ws.Range("H8").Value = 0.054
$ws.Range("I8").Value = 178
$ws.Range("J8").Value = 812

$ws.Range("B9").Value = 1.004
$ws.Range("C9").Value = 0.981
$ws.Range("D9").Value = -2.976
$ws.Range("E9").Value = 0.993
$ws.Range("F9").Value = 0.091
$ws.Range("G9").Value = 0.031
$ws.Range("H9").Value = 0.024
$ws.Range("I9").Value = 175
$ws.Range("J9").Value = 796

$ws.Range("B10").Value = 1.002
$ws.Range("C10").Value = 0.978
$ws.Range("D10").Value = -2.976
$ws.Range("E10").Value = 0.994
$ws.Range("F10").Value = 0.08599999999999999
$ws.Range("G10").Value = 0.033
$ws.Range("H10").Value = 0.024
$ws.Range("I10").Value = 192
$ws.Range("J10").Value = 871

$ws.Range("B11").Value = 1.002
$ws.Range("C11").Value = 0.98
$ws.Range("D11").Value = -2.975
$ws.Range("E11").Value = 0.992
$ws.Range("F11").Value = 0.081
$ws.Range("G11").Value = 0.033
$ws.Range("H11").Value = 0.025
$ws.Range("I11").Value = 162
$ws.Range("J11").Value = 730

$ws.Range("B12").Value = 1.001
$ws.Range("C12").Value = 0.984
$ws.Range("D12").Value = -2.971
$ws.Range("E12").Value = 0.987
$ws.Range("F12").Value = 0.077
$ws.Range("G12").Value = 0.036
$ws.Range("H12").Value = 0.029
$ws.Range("I12").Value = 142
$ws.Range("J12").Value = 633
